$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.014.42'
$ws.Range("E2").Value = '  -3.69%  '

$ws.Range("D3").Value = '2.360.10'
$ws.Range("E3").Value = '  -4.02%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '501.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.85%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.49%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.543'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.77%  '

$ws.Range("D9").Value = '2.363.16'
$ws.Range("E9").Value = '  -3.82%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0982'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.03%  '

$ws.Range("E11").Value = '  -0.01%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.78'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.324'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.21%  '

$ws.Range("D14").Value = '2.778.30'
$ws.Range("E14").Value = '  -3.89%  '

$ws.Range("D15").Value = '55.984.38'
$ws.Range("E15").Value = '  -3.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.06%  '

$ws.Range("E17").Value = '  -2.94%  '

$ws.Range("D18").Value = '2.372.29'
$ws.Range("E18").Value = '  -1.81%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.58%  '

$ws.Range("E20").Value = '  -2.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '306.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.99%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.52%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.67%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.998'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.364'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.146'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.21%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.20'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.05'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.92%  '

$ws.Range("D30").Value = '0.0₃0707'
$ws.Range("E30").Value = '  -4.25%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.64'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.29%  '

$ws.Range("E32").Value = '  +0.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.77'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.07'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.57'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.11%  '

$ws.Range("E37").Value = '  -4.86%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.69'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.38%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.22'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.41%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.803'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.73%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.37'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.52%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.34'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '128.10'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.98%  '

$ws.Range("E44").Value = '  -4.82%  '

$ws.Range("E45").Value = '  -1.96%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0892'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.78%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '238.01'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0480'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0205'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.93'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.948'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.25%  '
